# Simulated Wild Card round and logged it
# Add a new Running Back (K.Barner) row to the "RB" stats sheet with
# the round's stat line (all zeros for this simulated game).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("RB")

$row = 7

$ws.Cells.Item($row, 1).Value = "K.Barner"
$ws.Cells.Item($row, 2).Value = 0
$ws.Cells.Item($row, 3).Value = 0
$ws.Cells.Item($row, 4).Value = 0
$ws.Cells.Item($row, 5).Value = 0
$ws.Cells.Item($row, 6).Value = 0
$ws.Cells.Item($row, 7).Value = 0
$ws.Cells.Item($row, 8).Value = 0
$ws.Cells.Item($row, 9).Value = 0
$ws.Cells.Item($row, 10).Value = 0

$ws.Range("J8").Select()
